$d = $word.ActiveDocument

$d.Content.Find.Execute("82×50=4100", $true, $false, $false, $false, $false, $true, 1, $false, "17×80=1360", 2) | Out-Null
$d.Content.Find.Execute("55×24=1320", $true, $false, $false, $false, $false, $true, 1, $false, "35×62=2170", 2) | Out-Null
$d.Content.Find.Execute("26×59=1534", $true, $false, $false, $false, $false, $true, 1, $false, "74×72=5328", 2) | Out-Null
$d.Content.Find.Execute("30×34=1020", $true, $false, $false, $false, $false, $true, 1, $false, "51×39=1989", 2) | Out-Null
$d.Content.Find.Execute("45×32=1440", $true, $false, $false, $false, $false, $true, 1, $false, "43×57=2451", 2) | Out-Null
$d.Content.Find.Execute("72×60=4320", $true, $false, $false, $false, $false, $true, 1, $false, "40×92=3680", 2) | Out-Null
$d.Content.Find.Execute("29×79=2291", $true, $false, $false, $false, $false, $true, 1, $false, "18×94=1692", 2) | Out-Null
$d.Content.Find.Execute("47×57=2679", $true, $false, $false, $false, $false, $true, 1, $false, "94×80=7520", 2) | Out-Null
$d.Content.Find.Execute("39×84=3276", $true, $false, $false, $false, $false, $true, 1, $false, "34×93=3162", 2) | Out-Null
$d.Content.Find.Execute("43×44=1892", $true, $false, $false, $false, $false, $true, 1, $false, "84×83=6972", 2) | Out-Null
$d.Content.Find.Execute("99×48=4752", $true, $false, $false, $false, $false, $true, 1, $false, "61×37=2257", 2) | Out-Null
$d.Content.Find.Execute("52×36=1872", $true, $false, $false, $false, $false, $true, 1, $false, "66×12=792", 2) | Out-Null
$d.Content.Find.Execute("35×80=2800", $true, $false, $false, $false, $false, $true, 1, $false, "53×12=636", 2) | Out-Null
$d.Content.Find.Execute("88×38=3344", $true, $false, $false, $false, $false, $true, 1, $false, "76×84=6384", 2) | Out-Null
$d.Content.Find.Execute("21×44=924", $true, $false, $false, $false, $false, $true, 1, $false, "60×22=1320", 2) | Out-Null
$d.Content.Find.Execute("58×30=1740", $true, $false, $false, $false, $false, $true, 1, $false, "18×14=252", 2) | Out-Null
$d.Content.Find.Execute("65×97=6305", $true, $false, $false, $false, $false, $true, 1, $false, "84×39=3276", 2) | Out-Null
$d.Content.Find.Execute("44×69=3036", $true, $false, $false, $false, $false, $true, 1, $false, "38×29=1102", 2) | Out-Null
$d.Content.Find.Execute("59×56=3304", $true, $false, $false, $false, $false, $true, 1, $false, "94×30=2820", 2) | Out-Null
$d.Content.Find.Execute("82×26=2132", $true, $false, $false, $false, $false, $true, 1, $false, "32×22=704", 2) | Out-Null
$d.Content.Find.Execute("31×18=558", $true, $false, $false, $false, $false, $true, 1, $false, "93×51=4743", 2) | Out-Null
$d.Content.Find.Execute("79×52=4108", $true, $false, $false, $false, $false, $true, 1, $false, "32×48=1536", 2) | Out-Null
$d.Content.Find.Execute("26×74=1924", $true, $false, $false, $false, $false, $true, 1, $false, "83×86=7138", 2) | Out-Null
$d.Content.Find.Execute("23×93=2139", $true, $false, $false, $false, $false, $true, 1, $false, "82×22=1804", 2) | Out-Null
$d.Content.Find.Execute("92×65=5980", $true, $false, $false, $false, $false, $true, 1, $false, "60×15=900", 2) | Out-Null
